# Updates cryptos list prices / volume(1h) percentages.
# Values are written with a leading apostrophe (text-qualifier) so
# numeric-looking strings (e.g. "58.88") stay text instead of being
# coerced into numbers, then the cell style is reset to "Normal" so no
# stray quote-prefix / number-format style sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$updates = @(
    @{ Row = 2;  D = "37.040.24";  E = "  -0.23%  " },
    @{ Row = 3;  D = "2.055.58";   E = "  +0.02%  " },
    @{ Row = 4;  D = $null;        E = "  -0.23%  " },
    @{ Row = 5;  D = $null;        E = "  -1.42%  " },
    @{ Row = 6;  D = $null;        E = "  -2.07%  " },
    @{ Row = 7;  D = "58.88";      E = "  -1.61%  " },
    @{ Row = 8;  D = $null;        E = "  -0.03%  " },
    @{ Row = 9;  D = $null;        E = "  -3.13%  " },
    @{ Row = 10; D = "0.0773";     E = "  -2.85%  " },
    @{ Row = 11; D = $null;        E = "  +2.06%  " },
    @{ Row = 12; D = "15.47";      E = "  -4.11%  " },
    @{ Row = 13; D = $null;        E = "  +5.47%  " },
    @{ Row = 14; D = "2.355.66";   E = "  +0.01%  " },
    @{ Row = 15; D = $null;        E = "  -0.58%  " },
    @{ Row = 16; D = "2.032.33";   E = "  -1.15%  " },
    @{ Row = 17; D = "18.17";      E = "  -3.73%  " },
    @{ Row = 18; D = "36.994.30";  E = "  -0.37%  " },
    @{ Row = 19; D = "73.96";      E = "  -2.67%  " },
    @{ Row = 20; D = "0.0₃0890";   E = "  -2.02%  " },
    @{ Row = 21; D = $null;        E = "  -0.39%  " },
    @{ Row = 22; D = "238.11";     E = "  -0.31%  " },
    @{ Row = 23; D = $null;        E = "  -0.06%  " },
    @{ Row = 24; D = "2.45";       E = "  +1.04%  " },
    @{ Row = 25; D = "10.06";      E = "  +4.92%  " },
    @{ Row = 26; D = "169.65";     E = "  +0.02%  " },
    @{ Row = 27; D = $null;        E = "  -3.13%  " },
    @{ Row = 28; D = $null;        E = "  -1.16%  " },
    @{ Row = 29; D = "5.44";       E = "  +12.96%  " },
    @{ Row = 30; D = $null;        E = "  -2.16%  " },
    @{ Row = 31; D = $null;        E = "  -2.29%  " },
    @{ Row = 32; D = "4.68";       E = "  +2.70%  " },
    @{ Row = 33; D = "0.0618";     E = "  -2.47%  " },
    @{ Row = 34; D = $null;        E = "  -0.07%  " },
    @{ Row = 35; D = "2.31";       E = "  +3.33%  " },
    @{ Row = 36; D = "1.84";       E = "  +5.56%  " },
    @{ Row = 37; D = "0.0843";     E = "  -6.32%  " },
    @{ Row = 38; D = $null;        E = "  -0.93%  " },
    @{ Row = 39; D = $null;        E = "  +1.44%  " },
    @{ Row = 40; D = "3.06";       E = "  -1.26%  " },
    @{ Row = 41; D = $null;        E = "  -0.16%  " },
    @{ Row = 42; D = $null;        E = "  +1.37%  " },
    @{ Row = 43; D = $null;        E = "  -10.50%  " },
    @{ Row = 44; D = "97.66";      E = "  -0.33%  " },
    @{ Row = 45; D = "17.00";      E = "  -4.65%  " },
    @{ Row = 46; D = "1.302.11";   E = "  +0.37%  " },
    @{ Row = 47; D = $null;        E = "  -5.81%  " },
    @{ Row = 48; D = $null;        E = "  -0.49%  " },
    @{ Row = 49; D = "6.81";       E = "  -0.53%  " },
    @{ Row = 50; D = "2.243.65";   E = "  +0.15%  " },
    @{ Row = 51; D = "44.42";      E = "  +1.71%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextValue $ws.Cells.Item($u.Row, 4) $u.D
    }
    Set-TextValue $ws.Cells.Item($u.Row, 5) $u.E
}
